$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Delete the first 4 data rows (rows 2-5), shifting the rest up.
    $ws.Range("A2:A5").EntireRow.Delete()

    # Renumber column A as a 0-based sequential index for the remaining data rows.
    for ($i = 0; $i -le 14; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
